$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the old sub-header row (Hiver/Eté/Année/(MW)/(GWh) labels) - this
# shifts the data rows up by one.
$ws.Rows.Item(2).Delete()

# Rebuild row 1 as a full header row.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# New selection, matching the saved view state in the target file.
$ws.Range("A2:K2").Select()
